$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the last sample header in BCA rep1 (column O, row 1)
$ws.Range("O1").Value = "IDH1-M2-1-F"

# Update the corresponding measurement value (column O, row 2)
$ws.Range("O2").Value = 0.46419999000000001

# Update the active selection to match the saved view
$ws.Range("P6").Select()
